$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 46, shifting existing rows 46:98 down to 47:99
$ws.Rows(46).Insert()

# Populate the newly inserted row 46 with the new data record
$ws.Range("A46").Value = 10
$ws.Range("B46").Value = "Vega Modelo de Temuco"
$ws.Range("C46").Value = "La Araucanía"
$ws.Range("D46").Value = 44512
$ws.Range("E46").Value = 9
$ws.Range("F46").Value = 100112012
$ws.Range("G46").Value = "Espinaca"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 20
$ws.Range("K46").Value = 8000
$ws.Range("L46").Value = 8000
$ws.Range("M46").Value = 8000
$ws.Range("N46").Value = "$/docena de atados"
$ws.Range("O46").Value = "Región de La Araucanía"
$ws.Range("P46").Value = 2667
$ws.Range("Q46").Value = 3
$ws.Range("R46").Value = "Hortaliza"
